# Thu, Jun 18, 2020  5:05:43 PM
#
# Re-colour the deck's theme: the "Integral" colour scheme used by the
# slide master (ppt/theme/theme1.xml) is replaced with the stock
# "Office Theme" colour scheme (the scheme that, before this edit, only
# lived in ppt/theme/theme2.xml / the Notes Master).
#
# ThemeColorScheme.Item(n) slot order (matches <a:clrScheme> child order):
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1 6 accent2 7 accent3
#   8 accent4 9 accent5 10 accent6 11 hlink 12 folHlink
#
# RGB() packs a hex "RRGGBB" string the same way VBA's RGB(r,g,b) does
# (r + g*256 + b*65536) -- the value PowerPoint's ColorFormat.RGB expects.

function Hex2Rgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Target ("Office Theme") colours -- these are the values that used to
# sit in ppt/theme/theme2.xml before the swap.
$officeColors = @{
    1  = "000000"  # dk1
    2  = "FFFFFF"  # lt1
    3  = "44546A"  # dk2
    4  = "E7E6E6"  # lt2
    5  = "5B9BD5"  # accent1
    6  = "ED7D31"  # accent2
    7  = "A5A5A5"  # accent3
    8  = "FFC000"  # accent4
    9  = "4472C4"  # accent5
    10 = "70AD47"  # accent6
    11 = "0563C1"  # hlink
    12 = "954F72"  # folHlink
}

$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Item($i).RGB = Hex2Rgb $officeColors[$i]
}
